$d = $word.ActiveDocument

# Find the anchor: the paragraph that ends with "Janeiro: Editora Interciência , 2004."
$anchorRange = $d.Content.Duplicate
[void]$anchorRange.Find.Execute("Janeiro: Editora Interciência , 2004.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Find the end: the paragraph containing the copyright notice
$endRange = $d.Content.Duplicate
[void]$endRange.Find.Execute("Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Delete from right after the "Janeiro..." paragraph's mark (consuming the blank paragraph
# that follows it) through to the end of the copyright paragraph's mark. This removes the
# "Ver no Jupiter Salvar em pdf Salvar em docx" and "© 2020 ..." paragraphs along with the
# blank paragraph that preceded them, leaving a single blank paragraph before the page break.
$paraStart = $anchorRange.Paragraphs(1).Range.End
$paraEnd = $endRange.Paragraphs(1).Range.End

$deleteRange = $d.Range($paraStart, $paraEnd)
[void]$deleteRange.Delete()
